$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered save games) for rows 2-4, columns B:G
$ws.Range("B2").Value = 0.01293466051926884
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1135.458372537421

$ws.Range("B3").Value = 0.2917716402565462
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 261.3203778131603
$ws.Range("E3").Value = 1133.036916526867
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1396.304844062544

$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 22.3905356188092
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 27.82738278199502
